$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Source" column before the existing "Description" column (C),
# pushing Description to column D.
$ws.Columns.Item(3).Insert()

# New column header
$ws.Range("C1").Value = "Source"

# Fill every data row (2-85) with the source label
$ws.Range("C2:C85").Value = "VertebradosIbericos - Reptiles"

# Reflect the selection left by the author after filling the new column
$ws.Range("C2:C85").Select()

Write-Host "done"
